$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Articles" column (C) with the new reference lists, and
# set the matching row heights (wrap-text rows grow/shrink to fit).
$ws.Range("C2").Value = "5, 6, 9, 10, 13, 15-17, 19-21, 23, 25, 26, 29, 31-33,  36, 40, 44, 46-48, 50-52, 55, 60, 69, 71-75, 79, 80, 87, 88, 90, 91, 93"
$ws.Rows.Item(2).RowHeight = 108

$ws.Range("C3").Value = "6, 12, 13, 27, 30, 31, 50, 59, 62, 70, 71, 80, 86"
$ws.Rows.Item(3).RowHeight = 54

$ws.Range("C4").Value = "6, 8, 10, 13, 16, 17, 19, 21-23, 25, 28, 29, 31, 32,  36, 37, 43, 47, 48, 50, 54, 60, 70, 71, 74, 80, 85, 91"
$ws.Rows.Item(4).RowHeight = 90

$ws.Range("C5").Value = "6, 8-10, 13, 14, 16, 17, 19, 21-23, 29, 31, 33, 37, 43, 44, 46-48, 50, 60, 70, 72-76, 78, 80, 84, 86, 88, 90, 91"
$ws.Rows.Item(5).RowHeight = 90

$ws.Range("C6").Value = "7, 10, 13, 16, 17, 20, 23, 27, 32, 33, 50, 53-55, 70, 72-74, 76, 87"
$ws.Rows.Item(6).RowHeight = 72

$ws.Range("C7").Value = "9, 13, 14, 16, 17, 31, 48, 50, 62, 72, 73, 78, 84, 85, 91"
$ws.Rows.Item(7).RowHeight = 54

$ws.Range("C8").Value = "6-11, 13, 14, 16, 17, 19, 20, 21, 23, 25-33, 35, 37, 40, 43, 46-48, 50, 53, 59, 60, 66, 70-72, 74, 78-80, 88, 90, 91"
$ws.Rows.Item(8).RowHeight = 108

$ws.Range("C9").Value = "6, 51, 55, 79, 80, 84"
$ws.Rows.Item(9).RowHeight = 54

$ws.Range("C10").Value = "6, 9, 16, 29, 44, 47, 50, 71-73, 76, 80, 111"
$ws.Rows.Item(10).RowHeight = 54

$ws.Range("C11").Value = "111"
$ws.Rows.Item(11).RowHeight = 72

$ws.Range("C12").Value = "6, 8, 10, 13, 20, 52, 88, 90"
$ws.Rows.Item(12).RowHeight = 108

# Update the view: zoom to 78%, scroll so row 7 is at top, and select A9.
$win = $ws.Application.ActiveWindow
$win.Zoom = 78
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("A9").Select()
